$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("C2").Value = 101.0018
$ws.Range("C3").Value = 303.7467
$ws.Range("C4").Value = 521.3089
$ws.Range("C5").Value = 545.6003
$ws.Range("C6").Value = 1258.742
$ws.Range("C7").Value = 1146.043
$ws.Range("C8").Value = 2612.967
$ws.Range("C9").Value = 2074.043

$ws.Range("G15").Select()

$excel.CalculateFullRebuild()
